$d = $word.ActiveDocument

$d.Content.Find.Execute("-- Attaching core tidyverse packages ------------------------ tidyverse 2.0.0 --", $true, $false, $false, $false, $false, $true, 1, $false, "── Attaching core tidyverse packages ──────────────────────── tidyverse 2.0.0 ──", 2) | Out-Null
$d.Content.Find.Execute("v dplyr     1.1.4     v readr     2.1.6", $true, $false, $false, $false, $false, $true, 1, $false, "✔ dplyr     1.1.4     ✔ readr     2.1.6", 2) | Out-Null
$d.Content.Find.Execute("v forcats   1.0.1     v stringr   1.6.0", $true, $false, $false, $false, $false, $true, 1, $false, "✔ forcats   1.0.1     ✔ stringr   1.6.0", 2) | Out-Null
$d.Content.Find.Execute("v ggplot2   4.0.1     v tibble    3.3.1", $true, $false, $false, $false, $false, $true, 1, $false, "✔ ggplot2   4.0.1     ✔ tibble    3.3.1", 2) | Out-Null
$d.Content.Find.Execute("v lubridate 1.9.4     v tidyr     1.3.2", $true, $false, $false, $false, $false, $true, 1, $false, "✔ lubridate 1.9.4     ✔ tidyr     1.3.2", 2) | Out-Null
$d.Content.Find.Execute("v purrr     1.2.1     ", $true, $false, $false, $false, $false, $true, 1, $false, "✔ purrr     1.2.1     ", 2) | Out-Null
$d.Content.Find.Execute("-- Conflicts ------------------------------------------ tidyverse_conflicts() --", $true, $false, $false, $false, $false, $true, 1, $false, "── Conflicts ────────────────────────────────────────── tidyverse_conflicts() ──", 2) | Out-Null
$d.Content.Find.Execute("x dplyr::filter() masks stats::filter()", $true, $false, $false, $false, $false, $true, 1, $false, "✖ dplyr::filter() masks stats::filter()", 2) | Out-Null
$d.Content.Find.Execute("x dplyr::lag()    masks stats::lag()", $true, $false, $false, $false, $false, $true, 1, $false, "✖ dplyr::lag()    masks stats::lag()", 2) | Out-Null
$d.Content.Find.Execute("i Use the conflicted package (<http://conflicted.r-lib.org/>) to force all conflicts to become errors", $true, $false, $false, $false, $false, $true, 1, $false, "ℹ Use the conflicted package (<http://conflicted.r-lib.org/>) to force all conflicts to become errors", 2) | Out-Null
